# docs: update diagrams and images
#
# Third refactoring phase: `persons` -> `modules`. This particular deck
# (UndoRedoStartingStateListDiagram.pptx) only carries boilerplate
# slide-layout metadata changes that ride along with PowerPoint re-saving
# the file:
#   - mark the deck to remove personal info on save
#   - refresh the cached "datetimeFigureOut" field text (stamped the next
#     time the Insert > Header & Footer "Apply to All" date got re-cached)
#     on every slide layout's Date placeholder, 7/6/2018 -> 2/26/2019

$p = $ppt.ActivePresentation

# Best-effort: mirror PowerPoint's "Remove personal information from this
# file on save" option (File > Options > Trust Center > Privacy Options).
$p.RemovePersonalInformation = $true

# Walk every slide layout that hangs off the (single) slide master and
# update its Date placeholder's cached text. Access the layouts through
# Designs(1).SlideMaster.CustomLayouts, which enumerates them in the same
# order as the slideLayoutN.xml parts (slideLayout1.xml first, etc.).
$oldDate = "7/6/2018"
$newDate = "2/26/2019"

$slideMaster = $p.Designs.Item(1).SlideMaster
$customLayouts = $slideMaster.CustomLayouts

for ($li = 1; $li -le $customLayouts.Count; $li++) {
    $layout = $customLayouts.Item($li)
    $shapes = $layout.Shapes

    for ($si = 1; $si -le $shapes.Count; $si++) {
        $shp = $shapes.Item($si)

        $placeholderType = -1
        try {
            $placeholderType = $shp.PlaceholderFormat.Type
        } catch {
            $placeholderType = -1
        }

        # ppPlaceholderDate = 16
        if ($placeholderType -eq 16) {
            if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}
